# volunteers.xlsx prototype import edit
# - rename the "card" sheet to "volunteers"
# - clear the stray "Manik" entry that had leaked into the
#   volunteeredBefore column (H2), resetting it to 0
# - clear the notes value in K2 (keep its cell style)
# - move the active selection to I12 (cursor position after the edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "card" to "volunteers"
$ws.Name = "volunteers"

# 2. H2 was a shared string "Manik" (volunteeredBefore) -> becomes numeric 0
$ws.Range("H2").Value = 0

# 3. K2 (notes) had a value of 1 -> cleared, style (s="12") is preserved
$ws.Range("K2").ClearContents()

# 4. Update the selection/cursor shown when the workbook is reopened
$ws.Range("I12").Select()
